$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ID" header in B3 becomes "Idé"
$ws.Range("B3").Value = "Idé"

# Selection moves from A7 to D11
$ws.Range("D11").Select() | Out-Null

# Default column width widens slightly
$ws.StandardWidth = 11.55078125

# Header/footer font style name changes from "Regular" to "Normal"
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
